$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$s.Shapes.Item("副标题 2").Cut()
